$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room: insert 4 new rows after the current last data row (row 19),
#    copying the formatting of row 19 (default row height, correct per-column
#    styles) so the new rows 20-23 inherit the right cell styles (s="5" Date,
#    s="6" No, s="1" Change Title, s="2" Detail, s="8" File Name).
# ---------------------------------------------------------------------------
$ws.Rows("19:19").Copy()
$ws.Rows("20:23").Insert(-4121)   # xlShiftDown

# ---------------------------------------------------------------------------
# 2) Re-write rows 18-23 with their final content (this both replaces the
#    old rows 18-19 text and fills in the 4 new rows).
# ---------------------------------------------------------------------------

# Row 18: Install firebase by using 'npm install firebase --save'
$ws.Range("A18").Value = 44986
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = "Install firebase by using 'npm install firebase --save'"
$ws.Range("D18").Value = "npm install firebase --save"
$ws.Range("E18").Value = ""

# Row 19: Install firebase by using 'npm i @angular/fire --save'
$ws.Range("A19").Value = 44986
$ws.Range("B19").Value = 18
$ws.Range("C19").Value = "Install firebase by using 'npm i @angular/fire --save'"
$ws.Range("D19").Value = "npm i @angular/fire --save"
$ws.Range("E19").Value = ""

# Row 20: Add firebaseConfig into environment.ts, environment.prod.ts
$ws.Range("A20").Value = 44986
$ws.Range("B20").Value = 19
$ws.Range("C20").Value = "Add firebaseConfig into environment.ts, environment.prod.ts"
$ws.Range("D20").Value = "/*" + [char]10 + "    01 Mar 2023 wutthichair" + [char]10 + "      Add firebaseConfig" + [char]10 + "  */" + [char]10 + "  firebaseConfig : {" + [char]10 + "    apiKey: ""AIzaSyBNzq0vOMsBXMCrwD5rZKXmwNjIYWs0ZLg""," + [char]10 + "    authDomain: ""redbook-taa.firebaseapp.com""," + [char]10 + "    projectId: ""redbook-taa""," + [char]10 + "    storageBucket: ""redbook-taa.appspot.com""," + [char]10 + "    messagingSenderId: ""826040339698""," + [char]10 + "    appId: ""1:826040339698:web:e5781b8f67c7762402546b""" + [char]10 + "  }"
$ws.Range("E20").Value = "environment.ts" + [char]10 + "environment.prod.ts"

# Row 21: Import environment configuration and firebase modules into app.module.ts
$ws.Range("A21").Value = 44986
$ws.Range("B21").Value = 20
$ws.Range("C21").Value = "Import environment configuration and firebase modules into app.module.ts"
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = "app.module.ts"

# Row 22: Create filebaseUser interface by using 'ng g i ./@core/shared/interface/firebaseUser'
$ws.Range("A22").Value = 44986
$ws.Range("B22").Value = 20
$ws.Range("C22").Value = "Create filebaseUser interface by using 'ng g i ./@core/shared/interface/firebaseUser'"
$ws.Range("D22").Value = "ng g i ./@core/shared/interface/firebaseUser" + [char]10 + "add required properties that we will get from firebase after completed authenthication"
$ws.Range("E22").Value = "firebase-user.ts"

# Row 23: Create firbase Authentication Service by using 'ng g s ./@core/shared/services/firebaseAuthentication'
$ws.Range("A23").Value = 44986
$ws.Range("B23").Value = 20
$ws.Range("C23").Value = "Create firbase Authentication Service by using 'ng g s ./@core/shared/services/firebaseAuthentication'"
$ws.Range("D23").Value = "ng g s ./@core/shared/services/firebaseAuthentication"
$ws.Range("E23").Value = "firebase-authentication.service.ts"

# ---------------------------------------------------------------------------
# Row heights: rows 18, 19 and 21 go back to the sheet's default (auto) height;
# rows 20, 22 and 23 need the explicit wrapped-text heights from the source file.
# ---------------------------------------------------------------------------
$ws.Rows("18:18").EntireRow.AutoFit()
$ws.Rows("19:19").EntireRow.AutoFit()
$ws.Rows("21:21").EntireRow.AutoFit()
$ws.Rows("20:20").RowHeight = 180
$ws.Rows("22:22").RowHeight = 45
$ws.Rows("23:23").RowHeight = 30

# ---------------------------------------------------------------------------
# 3) Grow the table (ListObject) to cover the new rows.
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E23"))

# ---------------------------------------------------------------------------
# 4) Widen column E and move the view / selection the way the author left it.
# ---------------------------------------------------------------------------
$ws.Columns("E:E").ColumnWidth = 51.140625

$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("E23").Select()
